# Fix formatting issues introduced when scraping floating point numbers
# and tidy up a few provider-name strings that used a stray comma instead
# of a period as a separator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Provider / company name strings: replace a comma separator with a period ---
$ws.Range("E51").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"

$ws.Range("E61").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E90").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"

$ws.Range("E69").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F69").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E91").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F91").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"

# --- "Importe" column (H): these were scraped as European-style decimal
# strings ("1.234,56"). Reformat them as plain decimal text
# ("1234.56"): strip the "." thousands separator and turn the "," decimal
# separator into a ".". Column is kept as Text so Excel does not
# re-interpret the digits as a number and drop the trailing zeros. ---

$importeRange = $ws.Range("H2:H117")
$importeRange.NumberFormat = "@"

for ($r = 2; $r -le 117; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = $cell.Value2
    if ($old -ne $null) {
        $new = $old.ToString().Replace(".", "").Replace(",", ".")
        $cell.Value = $new
    }
}

# Restore the original (default) cell style now that the text is safely
# stored as a literal string - only the values needed to change, not the
# look/format of the column.
$importeRange.Style = "Normal"
